$d = $word.ActiveDocument
$p6 = $d.Paragraphs.Item(6)
$r = $p6.Range
$r.Delete()

$d.Content.Find.Execute("(Remember to ask during lab to give tutors access to the repository.)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Efolio Task 3.2:", 2)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para" $i ":" $p.Range.Text
}
